$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "2021" year column to the header row, matching the
# formatting already used by the preceding "2020" column (N4).
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 2021

# Add the corresponding 2021 data point to the data row, matching the
# formatting already used by the preceding data cell (N5).
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 1.5020015556876996

$excel.CutCopyMode = $false

# Reflect the new selection position left by the edit.
$ws.Range("Q5").Select()
